$d = $word.ActiveDocument

$replacements = @(
    @("78×74=5772", "54×51=2754"),
    @("74×55=4070", "70×88=6160"),
    @("90×52=4680", "41×58=2378"),
    @("30×26=780", "59×36=2124"),
    @("63×60=3780", "61×29=1769"),
    @("56×42=2352", "73×72=5256"),
    @("78×36=2808", "81×28=2268"),
    @("35×38=1330", "72×67=4824"),
    @("39×45=1755", "25×44=1100"),
    @("69×74=5106", "24×87=2088"),
    @("34×96=3264", "94×94=8836"),
    @("84×94=7896", "29×26=754"),
    @("68×48=3264", "22×61=1342"),
    @("47×38=1786", "38×98=3724"),
    @("17×56=952", "71×52=3692"),
    @("68×84=5712", "23×40=920"),
    @("58×62=3596", "93×68=6324"),
    @("40×34=1360", "58×69=4002"),
    @("20×96=1920", "50×31=1550"),
    @("18×91=1638", "12×43=516"),
    @("14×23=322", "89×41=3649"),
    @("26×36=936", "13×76=988"),
    @("76×41=3116", "45×19=855"),
    @("12×98=1176", "36×80=2880"),
    @("86×83=7138", "90×91=8190")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
